$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 11:52"

# Update data values (Casos totales, Casos activos, Recuperados, Muertes) per row
$ws.Range("B4").Value = 59784
$ws.Range("C4").Value = 35841
$ws.Range("D4").Value = 15895
$ws.Range("E4").Value = 8048

$ws.Range("B5").Value = 48158
$ws.Range("C5").Value = 13319
$ws.Range("D5").Value = 30031
$ws.Range("E5").Value = 4808

$ws.Range("B6").Value = 16589
$ws.Range("C6").Value = 6323
$ws.Range("D6").Value = 8556
$ws.Range("E6").Value = 1710

$ws.Range("B7").Value = 15706
$ws.Range("C7").Value = 5306
$ws.Range("D7").Value = 8004
$ws.Range("E7").Value = 2396

$ws.Range("B8").Value = 12564
$ws.Range("C8").Value = 9974
$ws.Range("D8").Value = 1335
$ws.Range("E8").Value = 1255

$ws.Range("B9").Value = 11913
$ws.Range("C9").Value = 5039
$ws.Range("D9").Value = 5706
$ws.Range("E9").Value = 1168

$ws.Range("B10").Value = 9328
$ws.Range("C10").Value = 1841
$ws.Range("D10").Value = 7075
$ws.Range("E10").Value = 412

$ws.Range("B14").Value = 5004
$ws.Range("C14").Value = 2081
$ws.Range("D14").Value = 2188
$ws.Range("E14").Value = 735

$ws.Range("B15").Value = 4759
$ws.Range("C15").Value = 1978
$ws.Range("D15").Value = 2349

$ws.Range("B17").Value = 3897
$ws.Range("C17").Value = 2064
$ws.Range("D17").Value = 1507
$ws.Range("E17").Value = 326

$ws.Range("B23").Value = 2751
$ws.Range("C23").Value = 1733
$ws.Range("D23").Value = 585
$ws.Range("E23").Value = 433

$ws.Range("B30").Value = 2255
$ws.Range("C30").Value = 780
$ws.Range("D30").Value = 1214
$ws.Range("E30").Value = 261

$ws.Range("B32").Value = 2187
$ws.Range("C32").Value = 1075
$ws.Range("D32").Value = 979
$ws.Range("E32").Value = 133

$ws.Range("B33").Value = 2115
$ws.Range("C33").Value = 1290
$ws.Range("D33").Value = 637
$ws.Range("E33").Value = 188

$ws.Range("B38").Value = 1475
$ws.Range("C38").Value = 1113
$ws.Range("D38").Value = 232
$ws.Range("E38").Value = 130

$ws.Range("B59").Value = 114
$ws.Range("C59").Value = 92
$ws.Range("D59").Value = 20
